$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''24.639.45'
$ws.Range('E2').Value = '  +11.25%  '
$ws.Range('D3').Value = '''1.680.80'
$ws.Range('E3').Value = '  +5.96%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '''305.59'
$ws.Range('E5').Value = '  +2.60%  '
$ws.Range('D6').Value = '''0.9954'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('D8').Value = '''0.3416'
$ws.Range('E8').Value = '  +2.49%  '
$ws.Range('D9').Value = '''48.20'
$ws.Range('E9').Value = '  +17.09%  '
$ws.Range('D10').Value = '''1.159'
$ws.Range('E10').Value = '  +4.34%  '
$ws.Range('E11').Value = '  +4.13%  '
$ws.Range('D12').Value = '''0.9988'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = '''6.081'
$ws.Range('E13').Value = '  +4.71%  '
$ws.Range('D14').Value = '''20.08'
$ws.Range('E14').Value = '  +3.85%  '
$ws.Range('D15').Value = '''6.682'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').Value = '''1.679.86'
$ws.Range('E16').Value = '  +6.13%  '
$ws.Range('E17').Value = '  +3.87%  '
$ws.Range('D18').Value = '''0.9953'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').Value = '''0.06641'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = '''80.52'
$ws.Range('E20').Value = '  +6.17%  '
$ws.Range('D21').Value = '''16.33'
$ws.Range('E21').Value = '  +3.51%  '
$ws.Range('D22').Value = '''6.057'
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('D23').Value = '''12.06'
$ws.Range('E23').Value = '  +4.40%  '
$ws.Range('D24').Value = '''24.601.61'
$ws.Range('E24').Value = '  +11.00%  '
$ws.Range('D25').Value = '''2.401'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = '''3.352'
$ws.Range('E26').Value = '  -3.21%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '''2.649'
$ws.Range('E27').Value = '  +6.41%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''152.32'
$ws.Range('E28').Value = '  +2.55%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''19.39'
$ws.Range('E29').Value = '  +1.55%  '
$ws.Range('B30').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C30').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D30').Value = '''1.866.77'
$ws.Range('E30').Value = '  +6.47%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = '''127.06'
$ws.Range('E31').Value = '  +4.78%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''6.209'
$ws.Range('E32').Value = '  +6.07%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = '''4.020'
$ws.Range('E33').Value = '  +2.47%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''0.9736'
$ws.Range('E34').Value = '  +6.39%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '''0.08395'
$ws.Range('E35').Value = '  +3.24%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''1.683'
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '''12.27'
$ws.Range('E37').Value = '  +5.18%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '''0.06370'
$ws.Range('E38').Value = '  +5.89%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '''5.280'
$ws.Range('E39').Value = '  +3.74%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''8.662'
$ws.Range('E40').Value = '  +4.91%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.02295'
$ws.Range('E41').Value = '  +5.60%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.232'
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '''0.2081'
$ws.Range('E43').Value = '  +5.48%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.6065'
$ws.Range('E44').Value = '  +5.43%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').Value = '''0.9934'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '''3.752'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''12.96'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = '''0.5844'
$ws.Range('E48').Value = '  +5.37%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''125.41'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '''1.996'
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.07210'
$ws.Range('E51').Value = '  +7.50%  '
